$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the price cells whose new values look numeric
# (e.g. '0.999', '171.55'), so Excel stores them as text -- matching the
# source workbook's inlineStr cells -- instead of auto-converting to numbers.
$textForceCells = @("D5", "D6", "D10", "D13", "D20", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D33", "D36", "D38", "D40", "D41", "D43", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Apply the updated values (new price/volume snapshot) ---
$ws.Range("D2").Value = '65.835.58'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '2.695.24'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("D5").Value = '609.67'
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").Value = '158.07'
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  +6.31%  '
$ws.Range("D10").Value = '6.02'
$ws.Range("E10").Value = '  +3.98%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '30.37'
$ws.Range("E13").Value = '  +4.39%  '
$ws.Range("E14").Value = '  +8.21%  '
$ws.Range("D15").Value = '3.178.11'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '65.692.06'
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '2.689.13'
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("E18").Value = '  +0.58%  '
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").Value = '359.48'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("E21").Value = '  +3.55%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = '70.75'
$ws.Range("E23").Value = '  +4.34%  '
$ws.Range("D24").Value = '9.84'
$ws.Range("E24").Value = '  +3.55%  '
$ws.Range("E25").Value = '  +14.01%  '
$ws.Range("D26").Value = '1.66'
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("D28").Value = '0.172'
$ws.Range("E28").Value = '  +5.32%  '
$ws.Range("D29").Value = '8.39'
$ws.Range("E29").Value = '  +3.64%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '2.21'
$ws.Range("E30").Value = '  +5.35%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '543.15'
$ws.Range("E31").Value = '  +5.91%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").Value = '1.80'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("E34").Value = '  +6.56%  '
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("D36").Value = '0.432'
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("D38").Value = '163.21'
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").Value = '171.55'
$ws.Range("E41").Value = '  +4.20%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").Value = '42.61'
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '0.0616'
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '23.64'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("D47").Value = '2.27'
$ws.Range("E47").Value = '  +4.54%  '
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").Value = '0.656'
$ws.Range("E49").Value = '  +1.64%  '
$ws.Range("D50").Value = '21.10'
$ws.Range("E50").Value = '  +9.47%  '
$ws.Range("D51").Value = '0.0993'
$ws.Range("E51").Value = '  +1.39%  '

# Restore the default cell style on the cells we force-formatted above,
# so only the values (not the formatting) differ from the original.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}